# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late" / heading / "Outstanding") columns
#   one place to the right (-> O/P/Q), and give the freshly inserted column
#   the same width as column M.
# - Make "Repayment schedule" the active sheet/tab, with K15 selected.
# - This naturally removes the "active" / tabSelected status from whichever
#   sheet previously held it ("Edit Repayment Schedule").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (pushes existing N:P -> O:Q).
$ws.Columns("N:N").Insert()

# Match the width Excel would have copied from the preceding column (M).
$ws.Columns("N:N").ColumnWidth = 10.2

# Make this sheet the active tab with the selection Excel left it in.
$ws.Activate()
$ws.Range("K15").Select()
